# Applies the cryptos list update described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.891.07"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.853.96"
$ws.Range("E3").Value = "  +1.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "697.66"
$ws.Range("E5").Value = "  -1.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.57"
$ws.Range("E6").Value = "  -0.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.853.14"
$ws.Range("E7").Value = "  +1.74%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("E10").Value = "  -0.48%  "
$ws.Range("E11").Value = "  -3.30%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("E13").Value = "  +0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.21"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.502.57"
$ws.Range("E15").Value = "  +1.69%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.810.51"
$ws.Range("E16").Value = "  +0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "70.979.23"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("E18").Value = "  -0.56%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.41"
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "498.68"
$ws.Range("E21").Value = "  +3.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.64"
$ws.Range("E22").Value = "  -5.00%  "
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.81"
$ws.Range("E24").Value = "  +1.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000148"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("E26").Value = "  +0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.20"
$ws.Range("E27").Value = "  -1.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.12"
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.01%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.14"
$ws.Range("E30").Value = "  +0.76%  "
$ws.Range("E31").Value = "  -0.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.26"
$ws.Range("E32").Value = "  -2.90%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.52"
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.810.45"
$ws.Range("E35").Value = "  +1.91%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.19"
$ws.Range("E36").Value = "  -0.27%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +0.64%  "
$ws.Range("E39").Value = "  +6.12%  "
$ws.Range("E40").Value = "  +8.65%  "
$ws.Range("E41").Value = "  -2.74%  "
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "163.89"
$ws.Range("E45").Value = "  +1.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000309"
$ws.Range("E46").Value = "  -5.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "48.99"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.38"
$ws.Range("E48").Value = "  -3.19%  "
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.48"
$ws.Range("E50").Value = "  -5.87%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.63"
$ws.Range("E51").Value = "  +1.09%  "
